# recordar_usuario.xlsx - ajustes a escenario de recordar usuario maximo de
# intentos fallidos

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# --- Row 2 ---
$ws.Range("B2").Value = "Alterno"
$ws.Range("F2").Value = "25130112"
$ws.Range("G2").Value = "USUCTDC3"
$ws.Range("H2").Value = 4321
$ws.Range("K2").Value = "Usuario o clave inválida. Inténtalo nuevamente"

# --- Row 3 ---
$ws.Range("B3").Value = "Acierto"
$ws.Range("F3").Value = "25130110"
$ws.Range("G3").Value = "USUCTDC1"
$ws.Range("H3").Value = 1234
$ws.Range("K3").Value = "El usuario ha sido enviado al correo electrónico."

# --- Row 4 ---
$ws.Range("F4").Value = "25130233"
$ws.Range("G4").Value = "USUCTDC3"
$ws.Range("H4").Value = 4321
$ws.Range("K4").Value = "¡Lo sentimos!"

# --- Row 5 ---
$ws.Range("F5").Value = "1989636240"
$ws.Range("G5").Value = "OSVPPRU16"

# --- Row 6 ---
$ws.Range("F6").Value = "25130114"
$ws.Range("G6").Value = "USUCTDC5"

# --- New (empty, underlined-style) row 11 ---
$ws.Range("G11").Font.Underline = $true
$ws.Range("K11").Font.Underline = $true

# --- View: scrolled/selected cell moved ---
$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Range("J10").Select()
